# Actualización automática 2025-11-27 09:30:10
# Applies updated sales figures across the three report sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" - per-client sales by product group
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("E32").Value = 74.29000000000001
$wsGrupo.Range("M32").Value = 4149.23

$wsGrupo.Range("L51").Value = 126.72

$wsGrupo.Range("M60").Value = 616.77

# Row 62 holds "X de 60" client-count summaries per product group.
$wsGrupo.Range("E62").Value = "3 de 60"
$wsGrupo.Range("L62").Value = "11 de 60"
$wsGrupo.Range("M62").Value = "14 de 60"

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL" - per-client monthly sales (noviembre column F)
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F32").Value = 7796.66
$wsMensual.Range("F51").Value = 4884.83
$wsMensual.Range("F60").Value = 1356.21
$wsMensual.Range("F62").Value = 62379.83

# ---------------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL" - budget vs. sales compliance by group
# ---------------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumpl.Range("D4").Value = 1596.77
$wsCumpl.Range("E4").Value = -807.39
$wsCumpl.Range("F4").Value = 2.022815374091059

$wsCumpl.Range("D11").Value = 8292.4
$wsCumpl.Range("E11").Value = 7855.6
$wsCumpl.Range("F11").Value = 0.5135248947238048

$wsCumpl.Range("D12").Value = 27144.99
$wsCumpl.Range("E12").Value = 23162.01
$wsCumpl.Range("F12").Value = 0.5395867374321666

$wsCumpl.Range("D14").Value = 65345.77
$wsCumpl.Range("E14").Value = 32516.11766749098
$wsCumpl.Range("F14").Value = 0.6677346161769103
